$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.546.00"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "1.813.87"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'308.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "'1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.4554"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").Value = "'0.3615"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.66%  "
$ws.Range("D9").Value = "'46.57"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").Value = "'0.07124"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "'0.8938"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.52%  "
$ws.Range("D12").Value = "'0.07807"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").Value = "'19.60"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("D14").Value = "1.839.36"
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "'5.312"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("D16").Value = "'6.366"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "'85.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "'1.006"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'0.000008555"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.86%  "
$ws.Range("D20").Value = "'1.005"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "26.566.24"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "'14.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.94%  "
$ws.Range("D23").Value = "'4.989"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").Value = "2.054.94"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "'10.59"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").Value = "'1.979"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.40%  "
$ws.Range("D27").Value = "'151.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("D28").Value = "'17.92"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.56%  "
$ws.Range("D29").Value = "'2.062"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.58%  "
$ws.Range("D30").Value = "'112.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").Value = "'4.901"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "'0.08729"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.79%  "
$ws.Range("D33").Value = "'3.129"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("D34").Value = "'2.855"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +14.73%  "
$ws.Range("D35").Value = "'4.470"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").Value = "'0.7280"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").Value = "'1.119"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").Value = "'1.004"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "'1.073"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("D40").Value = "'0.01945"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D41").Value = "'0.05115"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'2.887"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.67%  "
$ws.Range("D43").Value = "'0.5144"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.74%  "
$ws.Range("D44").Value = "'6.831"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "'0.1518"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("D46").Value = "'8.090"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4677"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.10%  "
$ws.Range("B48").Value = "PaxDollar"
$ws.Range("C48").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D48").Value = "'1.005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.14%  "
$ws.Range("D49").Value = "'9.989"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.21%  "
$ws.Range("D50").Value = "'101.61"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.36%  "
$ws.Range("D51").Value = "'1.584"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.47%  "
